# Update the table style ("Table_0" -> built-in "{33746DDD-0F52-4A17-B02F-89AC91D7260E}")
# on the table shown on slide 5 of the deck.
#
# Old style id: {6C53D076-43A7-42C2-A1CC-952E2D4C902B}
# New style id: {33746DDD-0F52-4A17-B02F-89AC91D7260E}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTable) {
        $sh.Table.ApplyStyle("{33746DDD-0F52-4A17-B02F-89AC91D7260E}")
    }
}
